$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: paragraph "Decided to use a stack to store empty cell grid
# references rather than a doubly linked list. ..."
#   - expand the opening clause to describe the two-array structure
#   - append a new trailing sentence about reusing the structure for
#     populated cells
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(26)
$r1 = $p1.Range

$old1 = "Decided to use a stack to store empty cell grid references rather than a "
$new1 = "Decided to use a structure with two arrays, and a variable to point to the top of the arrays similar to a stack to store empty cell grid references rather than a "
$r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

$p1 = $d.Paragraphs(26)
$r1 = $p1.Range
$r1.InsertAfter(" This can also be used to store references to populated cells when trying to remove values for the player_grid, however the obtained cell reference in the list will be random between 0 and the top.")

# ---------------------------------------------------------------------------
# Edit 2: after the "Wrote functions for fill_first_empty and solve. ..."
# paragraph, add four new paragraphs (with a blank separator paragraph
# before the last one) describing later progress. The superscript "rd" in
# "3rd" is applied last, once all paragraph/text insertions are done, so
# that the "current formatting" state used by later InsertParagraphAfter
# calls never picks up the superscript.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(28)
$r2 = $p2.Range

# blank paragraph
$r2.InsertParagraphAfter()

# "Memory allocation for grid moved into function. ..."
$d.Paragraphs(29).Range.InsertParagraphAfter()
$d.Paragraphs(30).Range.Text = "Memory allocation for grid moved into function. Had trouble with pointers when moving it in."

# "Wrote function to copy grid (... 3rd ...)." -- plain text for now
$d.Paragraphs(30).Range.InsertParagraphAfter()
$d.Paragraphs(31).Range.Text = "Wrote function to copy grid (so as to have a solution and player grid, potentially later a 3rd to tell which are original values, and which are user values)."

# blank paragraph
$d.Paragraphs(31).Range.InsertParagraphAfter()

# "Considered logic of removing numbers from grid to create player_grid. ..."
$d.Paragraphs(32).Range.InsertParagraphAfter()
$d.Paragraphs(33).Range.Text = "Considered logic of removing numbers from grid to create player_grid. Decided against trying to do it symmetrically or to start in corners/centre, as want to check if solvable for any matches except the number just removed. If solvable, means there is more than 1 unique solution. Due to needing to check this after each number removed, means there isn’t much point in starting with corners/centre. Will utilise stack like structure to store references of populated cells, however generate a random number to obtain between 0 and the top. Will also store frequency of clues remaining in a hash table. This can be used to ensure that at least SIZE-1 clues remain for all except 1 number at all times."

# Now apply superscript formatting to the "rd" of "3rd", scoped to that
# paragraph's range so no other "rd" occurrence in the document is touched.
$p31 = $d.Paragraphs(31)
$r31 = $p31.Range
$r31.Find.Execute("rd", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r31.Font.Superscript = $true

Write-Output "done"
